$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.291.68'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.880.05'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.685'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.40'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.354'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('E10').Value = '  +1.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0738'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.69%  '
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.48'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('D14').Value = '2.155.29'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.763'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.11%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.927.10'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.89'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('D18').Value = '35.325.11'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.50%  '
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.15'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.27'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('D31').Value = '4.128.44'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0585'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.840'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0721'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.43'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.90'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.33%  '
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('D45').Value = '1.303.76'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0795'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.54%  '
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.06'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.20'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.71%  '
